$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # PositiveExtra
$ws2 = $wb.Worksheets.Item(2)  # ProductList
$ws3 = $wb.Worksheets.Item(3)  # CAPTSetting
$ws4 = $wb.Worksheets.Item(4)  # SMAPTSetting

# --- Content updates: add product #18 (row 19) to every sheet for the UK24 setting ---

# Sheet1 PositiveExtra: A19 = 18, B19 = 12.5
$ws1.Range("A19").Value = 18
$ws1.Range("A19").NumberFormat = "0"
$ws1.Range("B19").Value = 12.5

# Sheet2 ProductList: A19 = 18, B19 = "UK24" (new shared string)
$ws2.Range("A19").Value = 18
$ws2.Range("A19").NumberFormat = "0"
$ws2.Range("B19").Value = "UK24"

# Sheet3 CAPTSetting: A19=18, B19="UK24", C19=50.5, D19=49, E19=100
$ws3.Range("A19").Value = 18
$ws3.Range("A19").NumberFormat = "0"
$ws3.Range("B19").Value = "UK24"
$ws3.Range("C19").Value = 50.5
$ws3.Range("D19").Value = 49
$ws3.Range("E19").Value = 100

# Sheet4 SMAPTSetting: A19=18, B19="UK24", C19=20, D19=19, E19=80
$ws4.Range("A19").Value = 18
$ws4.Range("A19").NumberFormat = "0"
$ws4.Range("B19").Value = "UK24"
$ws4.Range("C19").Value = 20
$ws4.Range("D19").Value = 19
$ws4.Range("E19").Value = 80

# --- View / selection updates ---
# Order matters: the last sheet activated below ends up the workbook's
# displayed (tabSelected / activeTab) sheet, matching the target state
# where "ProductList" (sheet index 2) is the active tab.

# Sheet1: scrolled so row 8 is at the top, final selection D15, no longer the active tab
$ws1.Activate() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 8
$aw.ScrollColumn = 1
$ws1.Range("D15").Select() | Out-Null

# Sheet3: scrolled so row 8 is at the top, final selection C19:E19
$ws3.Activate() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 8
$aw.ScrollColumn = 1
$ws3.Range("C19:E19").Select() | Out-Null

# Sheet4: scrolled so row 13 is at the top, final selection B20
$ws4.Activate() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 13
$aw.ScrollColumn = 1
$ws4.Range("B20").Select() | Out-Null

# Sheet2: activated last -> becomes the active/selected tab, scrolled so row 8
# is at the top, final selection B20
$ws2.Activate() | Out-Null
$aw = $excel.ActiveWindow
$aw.ScrollRow = 8
$aw.ScrollColumn = 1
$ws2.Range("B20").Select() | Out-Null

# --- Workbook window geometry ---
$aw = $excel.ActiveWindow
$aw.Left = 0
$aw.Top = 0
$aw.Width = 13068
$aw.Height = 4596
